# Auto-generated Excel COM-interop script
# Applies scraped market-data refresh values across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit's scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End
$ws.Range("H5").Value = 877.25
$ws.Range("I5").Value = 865.9091
$ws.Range("K5").Value = 865.9091
$ws.Range("M5").Value = -750.9091

# Row 9: Distill, My Heart
$ws.Range("H9").Value = 2432.0637
$ws.Range("I9").Value = 319.09525
$ws.Range("K9").Value = 319.09525
$ws.Range("M9").Value = -150.09525

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2557.5
$ws.Range("I138").Value = 877.4878
$ws.Range("J138").Value = 4088.1777
$ws.Range("K138").Value = 2632.4634
$ws.Range("L138").Value = 12264.5331
$ws.Range("M138").Value = 2507.5366
$ws.Range("N138").Value = -22544.5331

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6167.283
$ws.Range("I32").Value = 3819.2
$ws.Range("K32").Value = 3819.2
$ws.Range("M32").Value = -3532.2

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 50570.46
$ws.Range("I74").Value = 29831.531
$ws.Range("K74").Value = 29831.531
$ws.Range("M74").Value = -28957.531

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 50570.46
$ws.Range("I77").Value = 29831.531
$ws.Range("K77").Value = 149157.655
$ws.Range("M77").Value = -144789.655

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 30766.434
$ws.Range("I94").Value = 690.9545000000001
$ws.Range("J94").Value = 113474
$ws.Range("K94").Value = 690.9545000000001
$ws.Range("L94").Value = 113474
$ws.Range("M94").Value = -239.9545000000001
$ws.Range("N94").Value = -114376

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 1525.1666
$ws.Range("J105").Value = 2818.8
$ws.Range("L105").Value = 2818.8
$ws.Range("N105").Value = -6312.8

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3227.7292
$ws.Range("I58").Value = 2847.879
$ws.Range("K58").Value = 2847.879
$ws.Range("M58").Value = -2644.879

# Row 99: O Pine
$ws.Range("H99").Value = 4534.607
$ws.Range("I99").Value = 4621.15
$ws.Range("K99").Value = 4621.15
$ws.Range("M99").Value = -3123.15

# Row 107: Built to Last
$ws.Range("H107").Value = 29324.182
$ws.Range("I107").Value = 42173.227
$ws.Range("K107").Value = 42173.227
$ws.Range("M107").Value = -40253.227

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4534.607
$ws.Range("I126").Value = 4621.15
$ws.Range("K126").Value = 13863.45
$ws.Range("M126").Value = -11393.45

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 7582.353
$ws.Range("I132").Value = 2262.077
$ws.Range("J132").Value = 24873.25
$ws.Range("K132").Value = 6786.231000000001
$ws.Range("L132").Value = 74619.75
$ws.Range("M132").Value = -4256.231000000001
$ws.Range("N132").Value = -79679.75

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1472.585
$ws.Range("I134").Value = 1315.2188
$ws.Range("K134").Value = 3945.6564
$ws.Range("M134").Value = -1410.6564

# Row 136: Turali Quality
$ws.Range("H136").Value = 3227.7292
$ws.Range("I136").Value = 2847.879
$ws.Range("K136").Value = 8543.636999999999
$ws.Range("M136").Value = -5993.636999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 37: I Love Lamprey
$ws.Range("H37").Value = 114809.445
$ws.Range("J37").Value = 114809.445
$ws.Range("L37").Value = 344428.335
$ws.Range("N37").Value = -344652.335

# Row 62: Little Orphan Candy
$ws.Range("H62").Value = 8543.362999999999
$ws.Range("I62").Value = 4666
$ws.Range("J62").Value = 9997.375
$ws.Range("K62").Value = 13998
$ws.Range("L62").Value = 29992.125
$ws.Range("M62").Value = -13312
$ws.Range("N62").Value = -31364.125

# Row 65: Confections of Confession (L)
$ws.Range("H65").Value = 8543.362999999999
$ws.Range("I65").Value = 4666
$ws.Range("J65").Value = 9997.375
$ws.Range("K65").Value = 41994
$ws.Range("L65").Value = 89976.375
$ws.Range("M65").Value = -38562
$ws.Range("N65").Value = -96840.375

# Row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 7763.6787
$ws.Range("J138").Value = 9983.462
$ws.Range("L138").Value = 29950.386
$ws.Range("N138").Value = -40230.386

$ws = $wb.Worksheets.Item("GSM")
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 11813.263
$ws.Range("J113").Value = 6265.6665
$ws.Range("L113").Value = 6265.6665
$ws.Range("N113").Value = -10605.6665

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 62295.2
$ws.Range("I122").Value = 130040.71
$ws.Range("J122").Value = 3017.875
$ws.Range("K122").Value = 390122.13
$ws.Range("L122").Value = 9053.625
$ws.Range("M122").Value = -387672.13
$ws.Range("N122").Value = -13953.625

# Row 132: On Board for Lar
$ws.Range("H132").Value = 5236.222
$ws.Range("I132").Value = 5390
$ws.Range("K132").Value = 16170
$ws.Range("M132").Value = -13640

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 5610.8423
$ws.Range("I22").Value = 2334.6
$ws.Range("J22").Value = 6107.242
$ws.Range("K22").Value = 2334.6
$ws.Range("L22").Value = 6107.242
$ws.Range("M22").Value = -2039.6
$ws.Range("N22").Value = -6697.242

# Row 27: Fire and Hide
$ws.Range("H27").Value = 5610.8423
$ws.Range("I27").Value = 2334.6
$ws.Range("J27").Value = 6107.242
$ws.Range("K27").Value = 2334.6
$ws.Range("L27").Value = 6107.242
$ws.Range("M27").Value = -2227.6
$ws.Range("N27").Value = -6321.242

# Row 55: It''s Not a Job, It''s a Calling
$ws.Range("H55").Value = 2171.3635
$ws.Range("I55").Value = 770
$ws.Range("J55").Value = 3339.1667
$ws.Range("K55").Value = 770
$ws.Range("L55").Value = 3339.1667
$ws.Range("M55").Value = -597
$ws.Range("N55").Value = -3685.1667

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 3437.3845
$ws.Range("I61").Value = 3352.9473
$ws.Range("J61").Value = 3666.5715
$ws.Range("K61").Value = 3352.9473
$ws.Range("L61").Value = 3666.5715
$ws.Range("M61").Value = -3150.9473
$ws.Range("N61").Value = -4070.5715

# Row 113: Peace in Rest
$ws.Range("H113").Value = 3437.3845
$ws.Range("I113").Value = 3352.9473
$ws.Range("J113").Value = 3666.5715
$ws.Range("K113").Value = 3352.9473
$ws.Range("L113").Value = 3666.5715
$ws.Range("M113").Value = -1182.9473
$ws.Range("N113").Value = -8006.5715

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3201.7036
$ws.Range("I132").Value = 2679.7896
$ws.Range("J132").Value = 4441.25
$ws.Range("K132").Value = 8039.3688
$ws.Range("L132").Value = 13323.75
$ws.Range("M132").Value = -5509.3688
$ws.Range("N132").Value = -18383.75

# Row 136: Respect for Br''aax
$ws.Range("H136").Value = 4142.0625
$ws.Range("I136").Value = 3136.5557
$ws.Range("K136").Value = 9409.667099999999
$ws.Range("M136").Value = -6859.667099999999

# Row 137: Lending Artisans a Hand
$ws.Range("H137").Value = 62062.5
$ws.Range("J137").Value = 64928.57
$ws.Range("L137").Value = 64928.57
$ws.Range("N137").Value = -75128.57000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 9: A Taste for Dalmaticae
$ws.Range("H9").Value = 3803
$ws.Range("I9").Value = 3803
$ws.Range("K9").Value = 3803
$ws.Range("L9").Value = -3663

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 11099.6
$ws.Range("I62").Value = 3998
$ws.Range("J62").Value = 12875
$ws.Range("K62").Value = 3998
$ws.Range("L62").Value = 12875
$ws.Range("M62").Value = -3374
$ws.Range("N62").Value = -14123

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 11099.6
$ws.Range("I65").Value = 3998
$ws.Range("J65").Value = 12875
$ws.Range("K65").Value = 19990
$ws.Range("L65").Value = 64375
$ws.Range("M65").Value = -16870
$ws.Range("N65").Value = -70615

# Row 96: Skills on Display
$ws.Range("H96").Value = 5278.5
$ws.Range("I96").Value = 5278.5
$ws.Range("K96").Value = 5278.5
$ws.Range("M96").Value = -3905.5

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1946.8043
$ws.Range("I122").Value = 1705.7059
$ws.Range("J122").Value = 2629.9167
$ws.Range("K122").Value = 5117.1177
$ws.Range("L122").Value = 7889.750100000001
$ws.Range("M122").Value = -2667.1177
$ws.Range("N122").Value = -12789.7501

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 10254.143
$ws.Range("I132").Value = 7797
$ws.Range("J132").Value = 24997
$ws.Range("K132").Value = 23391
$ws.Range("L132").Value = 74991
$ws.Range("M132").Value = -20861
$ws.Range("N132").Value = -80051

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2584.5095
$ws.Range("I136").Value = 2260.6667
$ws.Range("K136").Value = 6782.000100000001
$ws.Range("M136").Value = -4232.000100000001
